# Append new timesheet entries (groups, permissions, TrainingTest javascript)
# to Sheet1, rows 49-53, right after the existing last row (48).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$hours = @{ 49 = 3; 50 = 5; 51 = 5; 52 = 5; 53 = 1.5 }
$dates = @{ 49 = 40966; 50 = 40968; 51 = 40972; 52 = 40973; 53 = 40974 }
$desc  = @{
    49 = "Refactoring, groups"
    50 = "Groups"
    51 = "Groups, Permissions & Managers refactoring"
    52 = "Groups, TrainingTest javascript"
    53 = "TrainingTest javascript & styles"
}

# Write the description (column B) cells in the same order the original
# author typed them in, so new shared-string table entries are interned
# in that exact order (52, 53, 51, 50, 49 -> new ids 44..48).
$textOrder = @(52, 53, 51, 50, 49)
foreach ($r in $textOrder) {
    $ws.Cells.Item($r, 2).Value = $desc[$r]
}

for ($r = 49; $r -le 53; $r++) {
    $ws.Cells.Item($r, 3).Value = $hours[$r]

    $ws.Cells.Item($r, 4).Value = $dates[$r]
    # Copy the date number-format from the row above instead of assigning
    # NumberFormat directly (which would mint a brand-new custom numFmt
    # entry instead of reusing the workbook's existing date style).
    $ws.Cells.Item($r - 1, 4).Copy() | Out-Null
    $ws.Cells.Item($r, 4).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

# Scroll the view down and select the last entered cell, like the author
# would after appending the new rows.
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B53").Select()
